$wb = $excel.ActiveWorkbook

$oldText = "February 03 2026 17.29.55 EST"
$newText = "February 03 2026 18.05.36 EST"

$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A2").Value = $wsAbout.Range("A2").Text.Replace($oldText, $newText)
$wsAbout.Range("A6").Value = $wsAbout.Range("A6").Text.Replace($oldText, $newText)

$wsData = $wb.Worksheets.Item("Boundaries and methane sources")
for ($row = 2; $row -le 10; $row++) {
    $cell = $wsData.Cells.Item($row, 19)
    $cell.Value = $cell.Text.Replace($oldText, $newText)
}
